$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "미트" entry (row 94, pointing to shared string index 92) is being
# removed from the nickname list. Deleting the whole row shifts every
# subsequent row up by one and Excel automatically recalculates the
# shared-string table (dropping the now-unused "미트" entry and
# renumbering references), matching the target diff exactly.
$ws.Rows.Item(94).Delete()
